$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 1122514.2
$ws.Range("I4").Value = 1262828.4
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 1262828.4
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = -1262714.4
$ws.Range("N4").Value = -229

$ws.Range("H18").Value = 10325.637
$ws.Range("I18").Value = 10950
$ws.Range("J18").Value = 4082
$ws.Range("K18").Value = 10950
$ws.Range("L18").Value = 4082
$ws.Range("M18").Value = -10666
$ws.Range("N18").Value = -4650

$ws.Range("H26").Value = 13166.333
$ws.Range("I26").Value = 14750
$ws.Range("J26").Value = 9999
$ws.Range("K26").Value = 14750
$ws.Range("L26").Value = 9999
$ws.Range("M26").Value = -14406
$ws.Range("N26").Value = -10687

$ws.Range("H33").Value = 956.4
$ws.Range("I33").Value = 956.4
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 956.4
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -727.4

$ws.Range("H106").Value = 9069
$ws.Range("I106").Value = 2547.2307
$ws.Range("J106").Value = 37330
$ws.Range("K106").Value = 2547.2307
$ws.Range("L106").Value = 37330
$ws.Range("M106").Value = -1916.2307
$ws.Range("N106").Value = -38592

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()

$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()

$ws.Range("H88").Value = 18182682
$ws.Range("I88").Value = 900.6
$ws.Range("J88").Value = 33334166
$ws.Range("K88").Value = 900.6
$ws.Range("L88").Value = 33334166
$ws.Range("M88").Value = -494.6
$ws.Range("N88").Value = -33334978

$ws.Range("H91").Value = 18182682
$ws.Range("I91").Value = 900.6
$ws.Range("J91").Value = 33334166
$ws.Range("K91").Value = 900.6
$ws.Range("L91").Value = 33334166
$ws.Range("M91").Value = 503.4
$ws.Range("N91").Value = -33336974

$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 3273
$ws.Range("I99").Value = 3273
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 3273
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -1775

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 10769485
$ws.Range("I3").Value = 11666858
$ws.Range("J3").Value = 1000
$ws.Range("K3").Value = 11666858
$ws.Range("L3").Value = 1000
$ws.Range("M3").Value = -11666745
$ws.Range("N3").Value = -1226

$ws.Range("H22").Value = 380
$ws.Range("I22").Value = 200
$ws.Range("J22").Value = 2000
$ws.Range("K22").Value = 200
$ws.Range("L22").Value = 2000
$ws.Range("M22").Value = 150
$ws.Range("N22").Value = -2700

$ws.Range("H31").Value = 42103.85
$ws.Range("I31").Value = 39936.668
$ws.Range("J31").Value = 50463
$ws.Range("K31").Value = 39936.668
$ws.Range("L31").Value = 50463
$ws.Range("M31").Value = -39641.668
$ws.Range("N31").Value = -51053

$ws.Range("H34").Value = 42103.85
$ws.Range("I34").Value = 39936.668
$ws.Range("J34").Value = 50463
$ws.Range("K34").Value = 39936.668
$ws.Range("L34").Value = 50463
$ws.Range("M34").Value = -39734.668
$ws.Range("N34").Value = -50867

$ws.Range("H97").Value = 31342.066
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 31342.066
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 31342.066
$ws.Range("N97").Value = -33324.066

$ws.Range("H106").Value = 36335.5
$ws.Range("I106").Value = 0
$ws.Range("J106").Value = 36335.5
$ws.Range("K106").Value = 0
$ws.Range("L106").Value = 36335.5
$ws.Range("M106").ClearContents()
$ws.Range("N106").Value = -38859.5

$ws.Range("H132").Value = 2704.158
$ws.Range("I132").Value = 2496
$ws.Range("J132").Value = 4473.5
$ws.Range("K132").Value = 7488
$ws.Range("L132").Value = 13420.5
$ws.Range("M132").Value = -4958
$ws.Range("N132").Value = -18480.5

$ws.Range("H134").Value = 4081.8667
$ws.Range("I134").Value = 3523.5454
$ws.Range("J134").Value = 5617.25
$ws.Range("K134").Value = 10570.6362
$ws.Range("L134").Value = 16851.75
$ws.Range("M134").Value = -8035.636200000001
$ws.Range("N134").Value = -21921.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 488.85715
$ws.Range("I2").Value = 629.6316
$ws.Range("J2").Value = 191.66667
$ws.Range("K2").Value = 3777.7896
$ws.Range("L2").Value = 1150.00002
$ws.Range("M2").Value = -3664.7896
$ws.Range("N2").Value = -1376.00002

$ws.Range("H7").Value = 53.75
$ws.Range("I7").Value = 67.111115
$ws.Range("J7").Value = 13.666667
$ws.Range("K7").Value = 201.333345
$ws.Range("L7").Value = 41.000001
$ws.Range("M7").Value = -89.33334500000001
$ws.Range("N7").Value = -265.000001

$ws.Range("H12").Value = 121.30769
$ws.Range("I12").Value = 79
$ws.Range("J12").Value = 140.11111
$ws.Range("K12").Value = 237
$ws.Range("L12").Value = 420.33333
$ws.Range("M12").Value = -64
$ws.Range("N12").Value = -766.3333299999999

$ws.Range("H63").Value = 800
$ws.Range("I63").Value = 800
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 2400
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -1651

$ws.Range("H66").Value = 800
$ws.Range("I66").Value = 800
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 7200
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -3456

$ws.Range("H75").Value = 4166.5
$ws.Range("I75").Value = 666.6667
$ws.Range("J75").Value = 7666.3335
$ws.Range("K75").Value = 2000.0001
$ws.Range("L75").Value = 22999.0005
$ws.Range("M75").Value = -1002.0001
$ws.Range("N75").Value = -24995.0005

$ws.Range("H78").Value = 4166.5
$ws.Range("I78").Value = 666.6667
$ws.Range("J78").Value = 7666.3335
$ws.Range("K78").Value = 6000.0003
$ws.Range("L78").Value = 68997.0015
$ws.Range("M78").Value = -1008.0003
$ws.Range("N78").Value = -78981.0015

$ws.Range("H103").Value = 587.4
$ws.Range("I103").Value = 462.5
$ws.Range("J103").Value = 774.75
$ws.Range("K103").Value = 1387.5
$ws.Range("L103").Value = 2324.25
$ws.Range("M103").Value = -508.5
$ws.Range("N103").Value = -4082.25

$ws.Range("H131").Value = 2156.6667
$ws.Range("I131").Value = 2107.5
$ws.Range("J131").Value = 2181.25
$ws.Range("K131").Value = 6322.5
$ws.Range("L131").Value = 6543.75
$ws.Range("M131").Value = -1282.5
$ws.Range("N131").Value = -16623.75

$ws.Range("H132").Value = 2552.111
$ws.Range("I132").Value = 1910
$ws.Range("J132").Value = 2799.077
$ws.Range("K132").Value = 17190
$ws.Range("L132").Value = 25191.693
$ws.Range("M132").Value = -14660
$ws.Range("N132").Value = -30251.693

$ws.Range("H137").Value = 3595.5
$ws.Range("I137").Value = 2903.3333
$ws.Range("J137").Value = 4287.6665
$ws.Range("K137").Value = 8709.999899999999
$ws.Range("L137").Value = 12862.9995
$ws.Range("M137").Value = -3609.999899999999
$ws.Range("N137").Value = -23062.9995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 46685.285
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 46685.285
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 46685.285
$ws.Range("M26").ClearContents()
$ws.Range("N26").Value = -47245.285

$ws.Range("H49").Value = 14370.2
$ws.Range("I49").Value = 14951
$ws.Range("J49").Value = 14225
$ws.Range("K49").Value = 14951
$ws.Range("L49").Value = 14225
$ws.Range("M49").Value = -14767
$ws.Range("N49").Value = -14593

$ws.Range("H50").Value = 46685.285
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 46685.285
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 46685.285
$ws.Range("M50").ClearContents()
$ws.Range("N50").Value = -47681.285

$ws.Range("H113").Value = 3539.75
$ws.Range("I113").Value = 3940.111
$ws.Range("J113").Value = 3025
$ws.Range("K113").Value = 3940.111
$ws.Range("L113").Value = 3025
$ws.Range("M113").Value = -1770.111
$ws.Range("N113").Value = -7365

$ws.Range("H130").Value = 68974
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 68974
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 68974
$ws.Range("N130").Value = -79014

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 100402.75
$ws.Range("I16").Value = 133813.33
$ws.Range("J16").Value = 171
$ws.Range("K16").Value = 133813.33
$ws.Range("L16").Value = 171
$ws.Range("M16").Value = -133643.33
$ws.Range("N16").Value = -511

$ws.Range("H130").Value = 0
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()

$ws.Range("H132").Value = 3985.3635
$ws.Range("I132").Value = 3541.2104
$ws.Range("J132").Value = 6798.3335
$ws.Range("K132").Value = 10623.6312
$ws.Range("L132").Value = 20395.0005
$ws.Range("M132").Value = -8093.6312
$ws.Range("N132").Value = -25455.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 8447416
$ws.Range("I3").Value = 12571149
$ws.Range("J3").Value = 199950
$ws.Range("K3").Value = 12571149
$ws.Range("L3").Value = 199950
$ws.Range("M3").Value = -12571035
$ws.Range("N3").Value = -200178

$ws.Range("H11").Value = 550
$ws.Range("I11").Value = 100
$ws.Range("J11").Value = 1000
$ws.Range("K11").Value = 100
$ws.Range("L11").Value = 1000
$ws.Range("M11").Value = 42
$ws.Range("N11").Value = -1284

$ws.Range("H94").Value = 30330
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 30330
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 30330
$ws.Range("N94").Value = -32132

$ws.Range("H132").Value = 8697.166999999999
$ws.Range("I132").Value = 10467.786
$ws.Range("J132").Value = 2500
$ws.Range("K132").Value = 31403.358
$ws.Range("L132").Value = 7500
$ws.Range("M132").Value = -28873.358
$ws.Range("N132").Value = -12560

$ws.Range("H136").Value = 3522.611
$ws.Range("I136").Value = 3263
$ws.Range("J136").Value = 5599.5
$ws.Range("K136").Value = 9789
$ws.Range("L136").Value = 16798.5
$ws.Range("M136").Value = -7239
$ws.Range("N136").Value = -21898.5

